$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(9).Insert()
$c7 = $ws.Cells.Item(7,9)
$src = $ws.Cells.Item(10,8)
$c7.Interior = $src.Interior()
Write-Output "done"
